# Split the "Similarly ... id's ..." sentence into several runs and add a
# comma after "Similarly", while changing "id's" to "ids" - and drop the
# stale grammar-check proofErr markers that bracketed the old wording.

$d = $word.ActiveDocument

# Locate the target paragraph robustly (rather than a hard-coded index).
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Similarly if products have same price*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the target paragraph"
}

$r = $target.Range

# Pull the paragraph's existing opening tag + <w:pPr> and a sample <w:rPr>
# straight out of the live document so the replacement keeps the same
# numbering / style / font settings without hard-coding them.
$owx = $r.WordOpenXML

$null = $owx -match '(?s)(<w:p[ >].*?</w:pPr>)'
$pPrefix = $matches[1]

$null = $owx -match '(?s)<w:r[ >].*?(<w:rPr>.*?</w:rPr>)'
$rPr = $matches[1]

$xml = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
          '<w:body>' +
            $pPrefix +
            '<w:r>' + $rPr + '<w:t>Similarly</w:t></w:r>' +
            '<w:r>' + $rPr + '<w:t>,</w:t></w:r>' +
            '<w:r>' + $rPr + '<w:t xml:space="preserve"> if products have same price, product </w:t></w:r>' +
            '<w:r>' + $rPr + '<w:t>ids</w:t></w:r>' +
            '<w:r>' + $rPr + '<w:t xml:space="preserve"> should be considered while displaying them in order</w:t></w:r>' +
          '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$r.InsertXML($xml)

Write-Host "Updated paragraph text: $($target.Range.Text)"
